$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "execution time" for the "web design" row (C4) to the new value
$ws.Range("C4").Value = "1h 02m"

# Move the active selection to C5, matching the saved cursor position
$ws.Range("C5").Select()
